$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1422613115598779
$ws.Range("D2").Value = 0.01416580112239885
$ws.Range("E2").Value = 0.07163464387278218
$ws.Range("F2").Value = 3.890905671548509
$ws.Range("G2").Value = 0.002606369570968986
$ws.Range("J2").Value = 0.1600343275239027
$ws.Range("K2").Value = 2.27125686690124
$ws.Range("M2").Value = 0.5893527142670223
$ws.Range("N2").Value = 3.101655713739717
# Row 3
$ws.Range("B3").Value = 0.1330881548836658
$ws.Range("D3").Value = 0.01414850831606262
$ws.Range("E3").Value = 0.0715947835129942
$ws.Range("F3").Value = 3.863764449682492
$ws.Range("G3").Value = 0.00261190013579189
$ws.Range("J3").Value = 0.1598344864677266
$ws.Range("K3").Value = 2.175152401469234
$ws.Range("M3").Value = 0.5722308427863396
$ws.Range("N3").Value = 3.113603573085427
# Row 4
$ws.Range("B4").Value = 0.1275300269704189
$ws.Range("D4").Value = 0.01414105888723149
$ws.Range("E4").Value = 0.07159654939040472
$ws.Range("F4").Value = 3.848949806636355
$ws.Range("G4").Value = 0.002615474439317477
$ws.Range("J4").Value = 0.1597749973499951
$ws.Range("K4").Value = 2.117605715954653
$ws.Range("M4").Value = 0.5620796055875488
$ws.Range("N4").Value = 3.121711541766402
# Row 5
$ws.Range("B5").Value = 0.1252837935248436
$ws.Range("D5").Value = 0.01413882344812833
$ws.Range("E5").Value = 0.07160387709704708
$ws.Range("F5").Value = 3.843376858906339
$ws.Range("G5").Value = 0.002616976041122554
$ws.Range("J5").Value = 0.1597666416301067
$ws.Range("K5").Value = 2.094520858690373
$ws.Range("M5").Value = 0.5580335793936513
$ws.Range("N5").Value = 3.125209382807199
# Row 6
$ws.Range("B6").Value = 0.1249119430415249
$ws.Range("D6").Value = 0.01413850069779965
$ws.Range("E6").Value = 0.0716054932618313
$ws.Range("F6").Value = 3.842479476279706
$ws.Range("G6").Value = 0.002617228105942448
$ws.Range("J6").Value = 0.1597662134258364
$ws.Range("K6").Value = 2.090709685048864
$ws.Range("M6").Value = 0.5573672096352453
$ws.Range("N6").Value = 3.125801890748946
# Row 7
$ws.Range("B7").Value = 0.1274996574550329
$ws.Range("D7").Value = 0.01414102549480845
$ws.Range("E7").Value = 0.07159662144511891
$ws.Range("F7").Value = 3.848872770113132
$ws.Range("G7").Value = 0.002615494507846865
$ws.Range("J7").Value = 0.1597748203492984
$ws.Range("K7").Value = 2.117292906293358
$ws.Range("M7").Value = 0.5620246726196285
$ws.Range("N7").Value = 3.121757930739662
# Row 8
$ws.Range("B8").Value = 0.1390830481026626
$ws.Range("D8").Value = 0.0141591829841099
$ws.Range("E8").Value = 0.07161545948112114
$ws.Range("F8").Value = 3.881162687617774
$ws.Range("G8").Value = 0.002608239552679685
$ws.Range("J8").Value = 0.159952296199954
$ws.Range("K8").Value = 2.237815568528845
$ws.Range("M8").Value = 0.5833738666184161
$ws.Range("N8").Value = 3.105614919562299
# Row 9
$ws.Range("B9").Value = 0.1623846637449304
$ws.Range("D9").Value = 0.01421979568968013
$ws.Range("E9").Value = 0.07186025615091118
$ws.Range("F9").Value = 3.959221532880662
$ws.Range("G9").Value = 0.002595421952896174
$ws.Range("J9").Value = 0.160802518780244
$ws.Range("K9").Value = 2.485852546382773
$ws.Range("M9").Value = 0.6281240166269129
$ws.Range("N9").Value = 3.080098899996287
# Row 10
$ws.Range("B10").Value = 0.1798608185796127
$ws.Range("D10").Value = 0.01427942354151046
$ws.Range("E10").Value = 0.07216647009961186
$ws.Range("F10").Value = 4.025648850029114
$ws.Range("G10").Value = 0.002586854104317937
$ws.Range("J10").Value = 0.1617344818778363
$ws.Range("K10").Value = 2.675366034033971
$ws.Range("M10").Value = 0.6627855793183173
$ws.Range("N10").Value = 3.065117672805002
# Row 11
$ws.Range("B11").Value = 0.1878884630005615
$ws.Range("D11").Value = 0.01430980233810608
$ws.Range("E11").Value = 0.07233315054234168
$ws.Range("F11").Value = 4.057860559406834
$ws.Range("G11").Value = 0.002583138651398327
$ws.Range("J11").Value = 0.162225463508733
$ws.Range("K11").Value = 2.763196506233101
$ws.Range("M11").Value = 0.6789471564470659
$ws.Range("N11").Value = 3.059124867210613
# Row 12
$ws.Range("B12").Value = 0.1909394369628927
$ws.Range("D12").Value = 0.01432177186942063
$ws.Range("E12").Value = 0.07240019889982996
$ws.Range("F12").Value = 4.070346510863146
$ws.Range("G12").Value = 0.002581757730208073
$ws.Range("J12").Value = 0.1624210427606343
$ws.Range("K12").Value = 2.796690982379573
$ws.Range("M12").Value = 0.6851241569895166
$ws.Range("N12").Value = 3.056974198788964
# Row 13
$ws.Range("B13").Value = 0.1902818637616832
$ws.Range("D13").Value = 0.01431917333502231
$ws.Range("E13").Value = 0.07238558416013952
$ws.Range("F13").Value = 4.067644602921405
$ws.Range("G13").Value = 0.002582053980477953
$ws.Range("J13").Value = 0.1623784915783091
$ws.Range("K13").Value = 2.789466863569373
$ws.Range("M13").Value = 0.6837912907178492
$ws.Range("N13").Value = 3.057432098406068
# Row 14
$ws.Range("B14").Value = 0.1881392469178138
$ws.Range("D14").Value = 0.01431077775620615
$ws.Range("E14").Value = 0.07233858792780801
$ws.Range("F14").Value = 4.058882004546291
$ws.Range("G14").Value = 0.002583024521428023
$ws.Range("J14").Value = 0.1622413603291122
$ws.Range("K14").Value = 2.765947395627336
$ws.Range("M14").Value = 0.6794541984978224
$ws.Range("N14").Value = 3.05894554940835
# Row 15
$ws.Range("B15").Value = 0.1868282724987864
$ws.Range("D15").Value = 0.01430569581584251
$ws.Range("E15").Value = 0.07231031297656187
$ws.Range("F15").Value = 4.053552218116096
$ws.Range("G15").Value = 0.002583622392251739
$ws.Range("J15").Value = 0.1621586213764559
$ws.Range("K15").Value = 2.751571705745562
$ws.Range("M15").Value = 0.676805032239514
$ws.Range("N15").Value = 3.059888050318634
# Row 16
$ws.Range("B16").Value = 0.1793377482551648
$ws.Range("D16").Value = 0.01427750350060286
$ws.Range("E16").Value = 0.07215612749851807
$ws.Range("F16").Value = 4.023583974224721
$ws.Range("G16").Value = 0.002587100569788708
$ws.Range("J16").Value = 0.1617037452444734
$ws.Range("K16").Value = 2.669658875732409
$ws.Range("M16").Value = 0.6617373353184561
$ws.Range("N16").Value = 3.065525899664365
# Row 17
$ws.Range("B17").Value = 0.1747623809600611
$ws.Range("D17").Value = 0.01426104018603347
$ws.Range("E17").Value = 0.07206854756016234
$ws.Range("F17").Value = 4.005711020478486
$ws.Range("G17").Value = 0.002589280856071611
$ws.Range("J17").Value = 0.1614418723825395
$ws.Range("K17").Value = 2.619824332191683
$ws.Range("M17").Value = 0.6525948967885284
$ws.Range("N17").Value = 3.069195442021041
# Row 18
$ws.Range("B18").Value = 0.1721380650389079
$ws.Range("D18").Value = 0.01425187729763344
$ws.Range("E18").Value = 0.0720207516646223
$ws.Range("F18").Value = 3.995618530810134
$ws.Range("G18").Value = 0.002590552048394412
$ws.Range("J18").Value = 0.1612975578857814
$ws.Range("K18").Value = 2.591313161999381
$ws.Range("M18").Value = 0.6473734597585619
$ws.Range("N18").Value = 3.071383431066153
# Row 19
$ws.Range("B19").Value = 0.1712507756747215
$ws.Range("D19").Value = 0.01424882759075174
$ws.Range("E19").Value = 0.07200501176066787
$ws.Range("F19").Value = 3.992233558718283
$ws.Range("G19").Value = 0.00259098540227921
$ws.Range("J19").Value = 0.1612497783561579
$ws.Range("K19").Value = 2.581685872615481
$ws.Range("M19").Value = 0.6456119247995957
$ws.Range("N19").Value = 3.072137520790307
# Row 20
$ws.Range("B20").Value = 0.1752486802484441
$ws.Range("D20").Value = 0.01426276104219149
$ws.Range("E20").Value = 0.07207760385962203
$ws.Range("F20").Value = 4.007594205061224
$ws.Range("G20").Value = 0.002589046987043711
$ws.Range("J20").Value = 0.1614690962386049
$ws.Range("K20").Value = 2.625113522976392
$ws.Range("M20").Value = 0.6535642879264216
$ws.Range("N20").Value = 3.068796802784391
# Row 21
$ws.Range("B21").Value = 0.1887682854275141
$ws.Range("D21").Value = 0.01431323111745009
$ws.Range("E21").Value = 0.07235228525666138
$ws.Range("F21").Value = 4.061447962481822
$ws.Range("G21").Value = 0.002582738744590968
$ws.Range("J21").Value = 0.1622813769256339
$ws.Range("K21").Value = 2.772849239542268
$ws.Range("M21").Value = 0.6807265590419576
$ws.Range("N21").Value = 3.058497787746049
# Row 22
$ws.Range("B22").Value = 0.1976686672283705
$ws.Range("D22").Value = 0.01434893010136573
$ws.Range("E22").Value = 0.07255470975429823
$ws.Range("F22").Value = 4.098324305859109
$ws.Range("G22").Value = 0.002578767657764615
$ws.Range("J22").Value = 0.1628685352988697
$ws.Range("K22").Value = 2.870773713717142
$ws.Range("M22").Value = 0.6988108493298739
$ws.Range("N22").Value = 3.052458736618476
# Row 23
$ws.Range("B23").Value = 0.1929124893687515
$ws.Range("D23").Value = 0.01432962920543535
$ws.Range("E23").Value = 0.07244457862745435
$ws.Range("F23").Value = 4.078488535087359
$ws.Range("G23").Value = 0.002580873267488423
$ws.Range("J23").Value = 0.1625500019206711
$ws.Range("K23").Value = 2.818383474936923
$ws.Range("M23").Value = 0.6891284191810314
$ws.Range("N23").Value = 3.055618430511245
# Row 24
$ws.Range("B24").Value = 0.175028805259501
$ws.Range("D24").Value = 0.01426198210184459
$ws.Range("E24").Value = 0.07207350154792103
$ws.Range("F24").Value = 4.006742247694461
$ws.Range("G24").Value = 0.002589152664004502
$ws.Range("J24").Value = 0.1614567688990718
$ws.Range("K24").Value = 2.622721845861861
$ws.Range("M24").Value = 0.6531259182197786
$ws.Range("N24").Value = 3.068976783666557
# Row 25
$ws.Range("B25").Value = 0.1560182801415522
$ws.Range("D25").Value = 0.01420073760440488
$ws.Range("E25").Value = 0.07177181721880999
$ws.Range("F25").Value = 3.936517233052967
$ws.Range("G25").Value = 0.002598739602500394
$ws.Range("J25").Value = 0.160518631032879
$ws.Range("K25").Value = 2.417484590550202
$ws.Range("M25").Value = 0.6157065423951948
$ws.Range("N25").Value = 3.086342036524883
